$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update coin Price (D) and Volume(1h) (E) columns with latest scraped values.
# NumberFormat '@' + reset to Normal style keeps numeric-looking price strings
# (e.g. '0.9989') stored as TEXT, matching the sheet's existing inlineStr cells,
# without leaving a stray number-format/quote-prefix style behind.
$cells = @(
  @{ Addr = "D2"; Val = '29.511.26' },
  @{ Addr = "E2"; Val = '  -0.70%  ' },
  @{ Addr = "D3"; Val = '1.849.63' },
  @{ Addr = "E3"; Val = '  -0.39%  ' },
  @{ Addr = "D4"; Val = '0.9989' },
  @{ Addr = "E4"; Val = '  +0.00%  ' },
  @{ Addr = "D5"; Val = '241.83' },
  @{ Addr = "E5"; Val = '  -1.15%  ' },
  @{ Addr = "D6"; Val = '0.6298' },
  @{ Addr = "E6"; Val = '  -1.71%  ' },
  @{ Addr = "D7"; Val = '0.9998' },
  @{ Addr = "E7"; Val = '  +0.03%  ' },
  @{ Addr = "D8"; Val = '47.85' },
  @{ Addr = "E8"; Val = '  +0.91%  ' },
  @{ Addr = "D9"; Val = '0.07566' },
  @{ Addr = "E9"; Val = '  +0.64%  ' },
  @{ Addr = "D10"; Val = '0.2976' },
  @{ Addr = "E10"; Val = '  -0.03%  ' },
  @{ Addr = "D11"; Val = '24.31' },
  @{ Addr = "E11"; Val = '  -0.59%  ' },
  @{ Addr = "D12"; Val = '0.07682' },
  @{ Addr = "E12"; Val = '  +0.24%  ' },
  @{ Addr = "D13"; Val = '1.892.23' },
  @{ Addr = "E13"; Val = '  +1.75%  ' },
  @{ Addr = "D14"; Val = '5.015' },
  @{ Addr = "E14"; Val = '  -0.40%  ' },
  @{ Addr = "E15"; Val = '  -0.91%  ' },
  @{ Addr = "D16"; Val = '83.75' },
  @{ Addr = "E16"; Val = '  -0.12%  ' },
  @{ Addr = "D17"; Val = '0.000009816' },
  @{ Addr = "E17"; Val = '  -0.10%  ' },
  @{ Addr = "D18"; Val = '2.123.23' },
  @{ Addr = "E18"; Val = '  +0.47%  ' },
  @{ Addr = "D19"; Val = '6.224' },
  @{ Addr = "E19"; Val = '  +2.25%  ' },
  @{ Addr = "D20"; Val = '29.555.70' },
  @{ Addr = "E20"; Val = '  -0.56%  ' },
  @{ Addr = "D21"; Val = '234.57' },
  @{ Addr = "E21"; Val = '  -0.74%  ' },
  @{ Addr = "E22"; Val = '  -1.35%  ' },
  @{ Addr = "D23"; Val = '0.9999' },
  @{ Addr = "E23"; Val = '  +0.01%  ' },
  @{ Addr = "D24"; Val = '7.613' },
  @{ Addr = "E24"; Val = '  +1.44%  ' },
  @{ Addr = "D25"; Val = '0.9999' },
  @{ Addr = "E25"; Val = '  -0.02%  ' },
  @{ Addr = "D26"; Val = '155.77' },
  @{ Addr = "E26"; Val = '  -1.91%  ' },
  @{ Addr = "D27"; Val = '0.1389' },
  @{ Addr = "E27"; Val = '  -2.19%  ' },
  @{ Addr = "D28"; Val = '8.437' },
  @{ Addr = "E28"; Val = '  -1.05%  ' },
  @{ Addr = "D29"; Val = '17.72' },
  @{ Addr = "D30"; Val = '1.482' },
  @{ Addr = "E30"; Val = '  -0.98%  ' },
  @{ Addr = "D31"; Val = '0.05837' },
  @{ Addr = "E31"; Val = '  -5.95%  ' },
  @{ Addr = "D32"; Val = '1.268' },
  @{ Addr = "E32"; Val = '  -1.30%  ' },
  @{ Addr = "E33"; Val = '  -1.32%  ' },
  @{ Addr = "D34"; Val = '4.036' },
  @{ Addr = "E34"; Val = '  -1.59%  ' },
  @{ Addr = "D35"; Val = '1.893' },
  @{ Addr = "E35"; Val = '  -0.39%  ' },
  @{ Addr = "E36"; Val = '  -0.02%  ' },
  @{ Addr = "D37"; Val = '0.7173' },
  @{ Addr = "E37"; Val = '  -1.58%  ' },
  @{ Addr = "D38"; Val = '2.589' },
  @{ Addr = "E38"; Val = '  -0.50%  ' },
  @{ Addr = "D39"; Val = '2.803' },
  @{ Addr = "E39"; Val = '  -0.88%  ' },
  @{ Addr = "D40"; Val = '1.235.32' },
  @{ Addr = "E40"; Val = '  +2.84%  ' },
  @{ Addr = "D41"; Val = '0.01781' },
  @{ Addr = "E41"; Val = '  -0.26%  ' },
  @{ Addr = "D42"; Val = '0.9114' },
  @{ Addr = "E42"; Val = '  -1.15%  ' },
  @{ Addr = "D43"; Val = '6.126' },
  @{ Addr = "E43"; Val = '  -1.68%  ' },
  @{ Addr = "D44"; Val = '2.032.31' },
  @{ Addr = "E44"; Val = '  +0.03%  ' },
  @{ Addr = "D45"; Val = '0.9993' },
  @{ Addr = "E45"; Val = '  -0.04%  ' },
  @{ Addr = "D46"; Val = '101.84' },
  @{ Addr = "E46"; Val = '  -0.12%  ' },
  @{ Addr = "D47"; Val = '67.49' },
  @{ Addr = "E47"; Val = '  +1.42%  ' },
  @{ Addr = "D48"; Val = '7.306' },
  @{ Addr = "E48"; Val = '  +9.45%  ' },
  @{ Addr = "D49"; Val = '9.169' },
  @{ Addr = "E49"; Val = '  -0.26%  ' },
  @{ Addr = "E50"; Val = '  -1.22%  ' },
  @{ Addr = "D51"; Val = '0.4030' }
)

foreach ($cell in $cells) {
  $rng = $ws.Range($cell.Addr)
  $rng.NumberFormat = "@"
  $rng.Value = $cell.Val
  $rng.Style = "Normal"
}
